{"js": "// Remove the trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line,\n// the \"\u00a9 2020 ...\" copyright/footer line, and the blank paragraph that sits\n// between them and the \"LOB1018: F\u00edsica I (Requisito fraco)\" requisite line \u2014\n// three whole paragraphs are deleted, leaving the requisite line followed\n// directly by the (already existing) trailing blank / page-break paragraphs.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Paragraphs identified unambiguously by their exact text content.\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Locate the \"LOB1018...\" requisite paragraph so we only remove the blank\n// paragraph that directly follows it (the document has several blank\n// paragraphs elsewhere that must stay untouched).\nlet lob1018Index = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"LOB1018: F\u00edsica I (Requisito fraco)\") {\n    lob1018Index = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (lob1018Index !== -1 && items[lob1018Index + 1] && items[lob1018Index + 1].text === \"\") {\n  toDelete.push(items[lob1018Index + 1]);\n}\nfor (let i = 0; i < items.length; i++) {\n  if (targetTexts.includes(items[i].text)) {\n    toDelete.push(items[i]);\n  }\n}\n\ntoDelete.forEach((para) => para.delete());\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line,\n# the \"\u00a9 2020 ...\" copyright/footer line, and the blank paragraph that sits\n# between them and the \"LOB1018: F\u00edsica I (Requisito fraco)\" requisite line \u2014\n# three whole paragraphs are deleted, leaving the requisite line followed\n# directly by the (already existing) trailing blank / page-break paragraphs.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaIndexByText($doc, $text) {\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        $t = $p.Range.Text -replace \"[\\r\\a\\f]+$\", \"\"\n        if ($t -eq $text) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# Locate the \"LOB1018...\" requisite paragraph so we only remove the blank\n# paragraph that directly follows it (the document has several blank\n# paragraphs elsewhere that must stay untouched).\n$lob1018Index = Get-ParaIndexByText $d \"LOB1018: F\u00edsica I (Requisito fraco)\"\n$blankIndex = $lob1018Index + 1\n\n$verIndex = Get-ParaIndexByText $d \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightIndex = Get-ParaIndexByText $d \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n# Delete from the highest paragraph index down to the lowest so earlier\n# deletions never invalidate the index of a paragraph still pending removal.\n$indices = @($blankIndex, $verIndex, $copyrightIndex) | Sort-Object -Descending\nforeach ($idx in $indices) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
